# manage order BOL update for staging
#
# Update the "BOL" sheet's order/tracking reference data:
#   - Order Id for the Parcel shipment (row 2) changes to the new order 59072066
#   - Carrier/tracking prefix "FCBTX" is replaced by "cev" for the LTL orders (rows 3-4)
# and leave the selection parked on the last edited cell, as Excel would.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOL")

$ws.Range("A2").Value = "59072066"
$ws.Range("C3").Value = "cev"
$ws.Range("C4").Value = "cev"

[void]$ws.Range("C4").Select()
